$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 23.67701833333334
$ws.Range("H2").Value = 71.03105500000001
$ws.Range("I2").Value = 0.1942294555737345
$ws.Range("J2").Value = 0.1942294555737345
$ws.Range("M2").Value = 4.233376
$ws.Range("N2").Value = 12.700128
$ws.Range("O2").Value = 0.1907241021355418
$ws.Range("P2").Value = 0.1907241021355418
$ws.Range("Q2").Value = 100.2337211638933
$ws.Range("R2").Value = 902.1034904750401
$ws.Range("S2").Value = 0.03704423852257561
$ws.Range("T2").Value = 0.03704423852257561

$ws.Range("G3").Value = 23.67701833333334
$ws.Range("H3").Value = 71.03105500000001
$ws.Range("I3").Value = 0.1942294555737345
$ws.Range("J3").Value = 0.1942294555737345
$ws.Range("M3").Value = 5.360490666666667
$ws.Range("N3").Value = 16.081472
$ws.Range("O3").Value = 0.2415034169905891
$ws.Range("P3").Value = 0.241503416990589
$ws.Range("Q3").Value = 126.9204357903289
$ws.Range("R3").Value = 1142.28392211296
$ws.Range("S3").Value = 0.04690707720127869
$ws.Range("T3").Value = 0.04690707720127869

$ws.Range("G4").Value = 23.67701833333334
$ws.Range("H4").Value = 71.03105500000001
$ws.Range("I4").Value = 0.1942294555737345
$ws.Range("J4").Value = 0.1942294555737345
$ws.Range("M4").Value = 0.147571
$ws.Range("N4").Value = 0.442713
$ws.Range("O4").Value = 0.006648440033732898
$ws.Range("P4").Value = 0.006648440033732897
$ws.Range("Q4").Value = 3.494041272468334
$ws.Range("R4").Value = 31.446371452215
$ws.Range("S4").Value = 0.001291322888166562
$ws.Range("T4").Value = 0.001291322888166562

$ws.Range("G5").Value = 23.67701833333334
$ws.Range("H5").Value = 71.03105500000001
$ws.Range("I5").Value = 0.1942294555737345
$ws.Range("J5").Value = 0.1942294555737345
$ws.Range("M5").Value = 10.554749
$ws.Range("N5").Value = 31.664247
$ws.Range("O5").Value = 0.4755176545364758
$ws.Range("P5").Value = 0.4755176545364758
$ws.Range("Q5").Value = 249.9049855767317
$ws.Range("R5").Value = 2249.144870190586
$ws.Range("S5").Value = 0.09235953515631885
$ws.Range("T5").Value = 0.09235953515631885

$ws.Range("G6").Value = 23.67701833333334
$ws.Range("H6").Value = 71.03105500000001
$ws.Range("I6").Value = 0.1942294555737345
$ws.Range("J6").Value = 0.1942294555737345
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2924906666666667
$ws.Range("N6").Value = 0.877472
$ws.Range("O6").Value = 0.01317743091637172
$ws.Range("P6").Value = 0.01317743091637172
$ws.Range("Q6").Value = 6.925306876995557
$ws.Range("R6").Value = 62.32776189296001
$ws.Range("S6").Value = 0.002559445232747377
$ws.Range("T6").Value = 0.002559445232747376

$ws.Range("G7").Value = 23.67701833333334
$ws.Range("H7").Value = 71.03105500000001
$ws.Range("I7").Value = 0.1942294555737345
$ws.Range("J7").Value = 0.1942294555737345
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.607657333333333
$ws.Range("N7").Value = 4.822972
$ws.Range("O7").Value = 0.07242895538728887
$ws.Range("P7").Value = 0.07242895538728887
$ws.Range("Q7").Value = 38.06453215505111
$ws.Range("R7").Value = 342.58078939546
$ws.Range("S7").Value = 0.01406783657264742
$ws.Range("T7").Value = 0.01406783657264742

$ws.Range("G8").Value = 14.89513733333333
$ws.Range("H8").Value = 44.685412
$ws.Range("I8").Value = 0.1221891360736233
$ws.Range("J8").Value = 0.1221891360736233
$ws.Range("M8").Value = 4.233376
$ws.Range("N8").Value = 12.700128
$ws.Range("O8").Value = 0.1907241021355418
$ws.Range("P8").Value = 0.1907241021355418
$ws.Range("Q8").Value = 63.05671690363733
$ws.Range("R8").Value = 567.510452132736
$ws.Range("S8").Value = 0.02330441326835934
$ws.Range("T8").Value = 0.02330441326835934

$ws.Range("G9").Value = 14.89513733333333
$ws.Range("H9").Value = 44.685412
$ws.Range("I9").Value = 0.1221891360736233
$ws.Range("J9").Value = 0.1221891360736233
$ws.Range("M9").Value = 5.360490666666667
$ws.Range("N9").Value = 16.081472
$ws.Range("O9").Value = 0.2415034169905891
$ws.Range("P9").Value = 0.241503416990589
$ws.Range("Q9").Value = 79.84524465405156
$ws.Range("R9").Value = 718.607201886464
$ws.Range("S9").Value = 0.02950909388090808
$ws.Range("T9").Value = 0.02950909388090808

$ws.Range("G10").Value = 14.89513733333333
$ws.Range("H10").Value = 44.685412
$ws.Range("I10").Value = 0.1221891360736233
$ws.Range("J10").Value = 0.1221891360736233
$ws.Range("M10").Value = 0.147571
$ws.Range("N10").Value = 0.442713
$ws.Range("O10").Value = 0.006648440033732898
$ws.Range("P10").Value = 0.006648440033732897
$ws.Range("Q10").Value = 2.198090311417333
$ws.Range("R10").Value = 19.782812802756
$ws.Range("S10").Value = 0.0008123671439591138
$ws.Range("T10").Value = 0.0008123671439591138

$ws.Range("G11").Value = 14.89513733333333
$ws.Range("H11").Value = 44.685412
$ws.Range("I11").Value = 0.1221891360736233
$ws.Range("J11").Value = 0.1221891360736233
$ws.Range("M11").Value = 10.554749
$ws.Range("N11").Value = 31.664247
$ws.Range("O11").Value = 0.4755176545364758
$ws.Range("P11").Value = 0.4755176545364758
$ws.Range("Q11").Value = 157.2144358738627
$ws.Range("R11").Value = 1414.929922864764
$ws.Range("S11").Value = 0.05810309139556764
$ws.Range("T11").Value = 0.05810309139556764

$ws.Range("G12").Value = 14.89513733333333
$ws.Range("H12").Value = 44.685412
$ws.Range("I12").Value = 0.1221891360736233
$ws.Range("J12").Value = 0.1221891360736233
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.2924906666666667
$ws.Range("N12").Value = 0.877472
$ws.Range("O12").Value = 0.01317743091637172
$ws.Range("P12").Value = 0.01317743091637172
$ws.Range("Q12").Value = 4.356688648718222
$ws.Range("R12").Value = 39.210197838464
$ws.Range("S12").Value = 0.001610138899341315
$ws.Range("T12").Value = 0.001610138899341315

$ws.Range("G13").Value = 14.89513733333333
$ws.Range("H13").Value = 44.685412
$ws.Range("I13").Value = 0.1221891360736233
$ws.Range("J13").Value = 0.1221891360736233
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.607657333333333
$ws.Range("N13").Value = 4.822972
$ws.Range("O13").Value = 0.07242895538728887
$ws.Range("P13").Value = 0.07242895538728887
$ws.Range("Q13").Value = 23.94627676494044
$ws.Range("R13").Value = 215.516490884464
$ws.Range("S13").Value = 0.008850031485487832
$ws.Range("T13").Value = 0.008850031485487833

$ws.Range("G14").Value = 83.33014933333334
$ws.Range("H14").Value = 249.990448
$ws.Range("I14").Value = 0.6835814083526421
$ws.Range("J14").Value = 0.6835814083526421
$ws.Range("M14").Value = 4.233376
$ws.Range("N14").Value = 12.700128
$ws.Range("O14").Value = 0.1907241021355418
$ws.Range("P14").Value = 0.1907241021355418
$ws.Range("Q14").Value = 352.7678542641493
$ws.Range("R14").Value = 3174.910688377344
$ws.Range("S14").Value = 0.1303754503446068
$ws.Range("T14").Value = 0.1303754503446068

$ws.Range("G15").Value = 83.33014933333334
$ws.Range("H15").Value = 249.990448
$ws.Range("I15").Value = 0.6835814083526421
$ws.Range("J15").Value = 0.6835814083526421
$ws.Range("M15").Value = 5.360490666666667
$ws.Range("N15").Value = 16.081472
$ws.Range("O15").Value = 0.2415034169905891
$ws.Range("P15").Value = 0.241503416990589
$ws.Range("Q15").Value = 446.690487753273
$ws.Range("R15").Value = 4020.214389779457
$ws.Range("S15").Value = 0.1650872459084023
$ws.Range("T15").Value = 0.1650872459084023

$ws.Range("G16").Value = 83.33014933333334
$ws.Range("H16").Value = 249.990448
$ws.Range("I16").Value = 0.6835814083526421
$ws.Range("J16").Value = 0.6835814083526421
$ws.Range("M16").Value = 0.147571
$ws.Range("N16").Value = 0.442713
$ws.Range("O16").Value = 0.006648440033732898
$ws.Range("P16").Value = 0.006648440033732897
$ws.Range("Q16").Value = 12.29711346726933
$ws.Range("R16").Value = 110.674021205424
$ws.Range("S16").Value = 0.004544750001607221
$ws.Range("T16").Value = 0.004544750001607221

$ws.Range("G17").Value = 83.33014933333334
$ws.Range("H17").Value = 249.990448
$ws.Range("I17").Value = 0.6835814083526421
$ws.Range("J17").Value = 0.6835814083526421
$ws.Range("M17").Value = 10.554749
$ws.Range("N17").Value = 31.664247
$ws.Range("O17").Value = 0.4755176545364758
$ws.Range("P17").Value = 0.4755176545364758
$ws.Range("Q17").Value = 879.5288103458508
$ws.Range("R17").Value = 7915.759293112657
$ws.Range("S17").Value = 0.3250550279845893
$ws.Range("T17").Value = 0.3250550279845892

$ws.Range("G18").Value = 83.33014933333334
$ws.Range("H18").Value = 249.990448
$ws.Range("I18").Value = 0.6835814083526421
$ws.Range("J18").Value = 0.6835814083526421
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 0.2924906666666667
$ws.Range("N18").Value = 0.877472
$ws.Range("O18").Value = 0.01317743091637172
$ws.Range("P18").Value = 0.01317743091637172
$ws.Range("Q18").Value = 24.37329093193956
$ws.Range("R18").Value = 219.359618387456
$ws.Range("S18").Value = 0.009007846784283027
$ws.Range("T18").Value = 0.009007846784283027

$ws.Range("G19").Value = 83.33014933333334
$ws.Range("H19").Value = 249.990448
$ws.Range("I19").Value = 0.6835814083526421
$ws.Range("J19").Value = 0.6835814083526421
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 1.607657333333333
$ws.Range("N19").Value = 4.822972
$ws.Range("O19").Value = 0.07242895538728887
$ws.Range("P19").Value = 0.07242895538728887
$ws.Range("Q19").Value = 133.9663256634951
$ws.Range("R19").Value = 1205.696930971456
$ws.Range("S19").Value = 0.04951108732915362
$ws.Range("T19").Value = 0.04951108732915362
